# Vitamin B3.docx: pick up new gram/mg data, fix "Hirse"/"Quinoa" rows,
# and tag the green-highlighted data lines as en-GB (adds proofErr markers
# for the now spell-checked food names).

$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [string]$findText,
        [string]$innerXml
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $findText"
    }
    $para = $rng.Paragraphs(1)
    $target = $para.Range

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$target.InsertXML($pkg)
}


# --- "Banane 100 g" -> "Banane" + " 100 g" (spell-checked, en-GB), plus
#     en-GB tagging of the rest of the line ---------------------------------
$banane = '<w:p w14:paraId="16C049F3" w14:textId="048F1C76" w:rsidR="00D03403" w:rsidRPr="0007024F" w:rsidRDefault="00D03403" w:rsidP="00D03403">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t>Banane</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> 100 g</w:t></w:r>' +
    '<w:r w:rsidR="00C00C17" w:rsidRPr="0007024F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00251C66" w:rsidRPr="0007024F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t>0,95 mg</w:t></w:r>' +
    '<w:r w:rsidRPr="0007024F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "Banane 100 g" $banane

# --- "Honigmelone 100 g" -> "Honigmelone" + " 100 g" (spell-checked,
#     en-GB), plus en-GB tagging of the rest of the line -------------------
$honigmelone = '<w:p w14:paraId="69E1B2B6" w14:textId="71022277" w:rsidR="00D03403" w:rsidRPr="0007024F" w:rsidRDefault="00D03403" w:rsidP="00D03403">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t>Honigmelone</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> 100 g</w:t></w:r>' +
    '<w:r w:rsidR="00C00C17" w:rsidRPr="0007024F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> 0,78 mg</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "Honigmelone 100 g" $honigmelone

# --- "Hirse 400g 5,7 mg" -> split so the corrected gram amount reads
#     cleanly as "Hirse 400" + " " + "g 5,7 mg" -----------------------------
$hirse = '<w:p w14:paraId="222A97DD" w14:textId="6FB50E58" w:rsidR="006D20F4" w:rsidRPr="006F7A51" w:rsidRDefault="006D20F4" w:rsidP="006D20F4">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="006D20F4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t>Hirse 400</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/><w:lang w:val="en-GB"/></w:rPr><w:t>g 5,7 mg</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "Hirse 400g 5,7 mg" $hirse

# --- "Quinoa ?" (unknown, red) -> "Quinoa 390 g 4,18 mg" (known, green) --
$quinoa = '<w:p w14:paraId="56C6B015" w14:textId="77777777" w:rsidR="00D03403" w:rsidRPr="00074C9F" w:rsidRDefault="00D03403" w:rsidP="00D03403">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Quinoa </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>390 g 4,18 mg</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml "Quinoa ?" $quinoa

Write-Output "applied 4 paragraph updates"
